# Auto-generated PowerShell Excel COM-interop script
# Applies updated TPM values to Lrpap1-Sorl1.xlsx sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.313941333333333
$ws.Range("H2").Value = 15.941824
$ws.Range("I2").Value = 0.176869630377001
$ws.Range("J2").Value = 0.176869630377001
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.033584
$ws.Range("N2").Value = 0.100752
$ws.Range("O2").Value = 0.002172419590320632
$ws.Range("P2").Value = 0.002172419590320632
$ws.Range("Q2").Value = 0.1784634057386666
$ws.Range("R2").Value = 1.606170651648
$ws.Range("S2").Value = 0.000384235049963766
$ws.Range("T2").Value = 0.0003842350499637661
$ws.Range("G3").Value = 5.313941333333333
$ws.Range("H3").Value = 15.941824
$ws.Range("I3").Value = 0.176869630377001
$ws.Range("J3").Value = 0.176869630377001
$ws.Range("O3").Value = 0.8824690642271135
$ws.Range("P3").Value = 0.8824690642271135
$ws.Range("Q3").Value = 72.49448281661867
$ws.Range("R3").Value = 652.450345349568
$ws.Range("S3").Value = 0.1560819772089875
$ws.Range("T3").Value = 0.1560819772089875
$ws.Range("G4").Value = 5.313941333333333
$ws.Range("H4").Value = 15.941824
$ws.Range("I4").Value = 0.176869630377001
$ws.Range("J4").Value = 0.176869630377001
$ws.Range("O4").Value = 0.115358516182566
$ws.Range("P4").Value = 0.115358516182566
$ws.Range("Q4").Value = 9.47665624570311
$ws.Range("R4").Value = 85.28990621132799
$ws.Range("S4").Value = 0.02040341811804973
$ws.Range("T4").Value = 0.02040341811804973
$ws.Range("I5").Value = 0.5461014638447835
$ws.Range("J5").Value = 0.5461014638447835
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.033584
$ws.Range("N5").Value = 0.100752
$ws.Range("O5").Value = 0.002172419590320632
$ws.Range("P5").Value = 0.002172419590320632
$ws.Range("Q5").Value = 0.551022393776
$ws.Range("R5").Value = 4.959201543984
$ws.Range("S5").Value = 0.001186361518359182
$ws.Range("T5").Value = 0.001186361518359182
$ws.Range("I6").Value = 0.5461014638447835
$ws.Range("J6").Value = 0.5461014638447835
$ws.Range("O6").Value = 0.8824690642271135
$ws.Range("P6").Value = 0.8824690642271135
$ws.Range("S6").Value = 0.4819176477721629
$ws.Range("T6").Value = 0.4819176477721629
$ws.Range("I7").Value = 0.5461014638447835
$ws.Range("J7").Value = 0.5461014638447835
$ws.Range("O7").Value = 0.115358516182566
$ws.Range("P7").Value = 0.115358516182566
$ws.Range("S7").Value = 0.06299745455426142
$ws.Range("T7").Value = 0.06299745455426142
$ws.Range("I8").Value = 0.2770289057782155
$ws.Range("J8").Value = 0.2770289057782155
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.033584
$ws.Range("N8").Value = 0.100752
$ws.Range("O8").Value = 0.002172419590320632
$ws.Range("P8").Value = 0.002172419590320632
$ws.Range("Q8").Value = 0.2795252181386666
$ws.Range("R8").Value = 2.515726963248
$ws.Range("S8").Value = 0.000601823021997684
$ws.Range("T8").Value = 0.0006018230219976841
$ws.Range("I9").Value = 0.2770289057782155
$ws.Range("J9").Value = 0.2770289057782155
$ws.Range("O9").Value = 0.8824690642271135
$ws.Range("P9").Value = 0.8824690642271135
$ws.Range("S9").Value = 0.2444694392459631
$ws.Range("T9").Value = 0.2444694392459631
$ws.Range("I10").Value = 0.2770289057782155
$ws.Range("J10").Value = 0.2770289057782155
$ws.Range("O10").Value = 0.115358516182566
$ws.Range("P10").Value = 0.115358516182566
$ws.Range("S10").Value = 0.03195764351025482
$ws.Range("T10").Value = 0.03195764351025482
